# Adding another notebooks for image classification
# (workbook data tweak: a few cells that were showing as mis-parsed dates
# need to hold their real decimal-gram text values instead)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a TEXT value into a cell while preserving that cell's
# existing number format / style (a plain Range.Value assignment of a
# numeric-looking string like "1.8" gets auto-coerced to a number, and
# forcing text via NumberFormat="@" or a leading apostrophe stamps a
# brand-new "Text" style onto the cell). Instead we compute the text in
# a scratch cell via a formula (formula-string results never need a
# special number format), copy it, and Paste-Special "Values only" into
# the destination - Excel's Values-only paste keeps the destination's
# own formatting and just swaps in the pasted (text) content.
function Set-TextKeepFormat {
    param($Sheet, [string]$Address, [string]$Text)

    $scratch = $Sheet.Range("ZZ1000")
    $scratch.Formula = '="' + $Text + '"'
    $scratch.Copy() | Out-Null
    $Sheet.Range($Address).PasteSpecial(-4163) | Out-Null
    $scratch.Clear() | Out-Null
}

Set-TextKeepFormat $ws "J2" "1.8"
Set-TextKeepFormat $ws "K2" "4.8"
Set-TextKeepFormat $ws "I6" "2.4"

# Leave the cursor where the author last left it.
$ws.Range("O8").Select() | Out-Null
